$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.247.25"
$ws.Range("E2").Value = "'  +0.35%  "
$ws.Range("D3").Value = "'1.861.79"
$ws.Range("E3").Value = "'  +0.65%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'0.7017"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("D6").Value = "'237.68"
$ws.Range("E6").Value = "'  -0.20%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'0.08130"
$ws.Range("E8").Value = "'  +8.93%  "
$ws.Range("D9").Value = "'0.3025"
$ws.Range("E9").Value = "'  -0.33%  "
$ws.Range("D10").Value = "'23.21"
$ws.Range("E10").Value = "'  -0.45%  "
$ws.Range("D11").Value = "'0.08143"
$ws.Range("E11").Value = "'  +0.18%  "
$ws.Range("D12").Value = "'1.840.93"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("D13").Value = "'5.162"
$ws.Range("E13").Value = "'  -1.02%  "
$ws.Range("D14").Value = "'0.7058"
$ws.Range("E14").Value = "'  -2.55%  "
$ws.Range("D15").Value = "'89.07"
$ws.Range("E15").Value = "'  +0.47%  "
$ws.Range("D16").Value = "'29.256.32"
$ws.Range("E16").Value = "'  +0.40%  "
$ws.Range("D17").Value = "'5.769"
$ws.Range("E17").Value = "'  +0.22%  "
$ws.Range("D18").Value = "'0.000007845"
$ws.Range("E18").Value = "'  +2.85%  "
$ws.Range("D19").Value = "'13.33"
$ws.Range("E19").Value = "'  +2.15%  "
$ws.Range("D20").Value = "'235.50"
$ws.Range("E20").Value = "'  -0.42%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "'  +0.04%  "
$ws.Range("D22").Value = "'2.109.34"
$ws.Range("E22").Value = "'  +0.73%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  +0.05%  "
$ws.Range("D24").Value = "'7.391"
$ws.Range("E24").Value = "'  -1.99%  "
$ws.Range("D25").Value = "'161.22"
$ws.Range("E25").Value = "'  +0.23%  "
$ws.Range("E26").Value = "'  -0.40%  "
$ws.Range("D27").Value = "'0.1438"
$ws.Range("E27").Value = "'  -0.73%  "
$ws.Range("D28").Value = "'18.07"
$ws.Range("E28").Value = "'  +0.27%  "
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("D30").Value = "'1.434"
$ws.Range("E30").Value = "'  +2.54%  "
$ws.Range("D31").Value = "'1.481"
$ws.Range("E31").Value = "'  -0.66%  "
$ws.Range("D32").Value = "'4.386"
$ws.Range("E32").Value = "'  -2.85%  "
$ws.Range("D33").Value = "'4.048"
$ws.Range("E33").Value = "'  +2.27%  "
$ws.Range("E34").Value = "'  +0.73%  "
$ws.Range("D35").Value = "'1.167"
$ws.Range("E35").Value = "'  -1.40%  "
$ws.Range("D36").Value = "'0.7049"
$ws.Range("E36").Value = "'  +1.04%  "
$ws.Range("D37").Value = "'0.9976"
$ws.Range("E37").Value = "'  -3.45%  "
$ws.Range("D38").Value = "'2.675"
$ws.Range("E38").Value = "'  +0.57%  "
$ws.Range("D39").Value = "'0.01840"
$ws.Range("E39").Value = "'  -1.20%  "
$ws.Range("D40").Value = "'2.731"
$ws.Range("E40").Value = "'  +2.06%  "
$ws.Range("D41").Value = "'0.9197"
$ws.Range("E41").Value = "'  -1.99%  "
$ws.Range("D42").Value = "'1.131.92"
$ws.Range("E42").Value = "'  +4.94%  "
$ws.Range("D43").Value = "'0.4268"
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("D44").Value = "'5.869"
$ws.Range("E44").Value = "'  -1.80%  "
$ws.Range("D45").Value = "'70.11"
$ws.Range("E45").Value = "'  +0.59%  "
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("D47").Value = "'102.69"
$ws.Range("E47").Value = "'  +0.74%  "
$ws.Range("D48").Value = "'1.763"
$ws.Range("E48").Value = "'  +1.47%  "
$ws.Range("D49").Value = "'2.003.21"
$ws.Range("E49").Value = "'  +0.83%  "
$ws.Range("D50").Value = "'9.167"
$ws.Range("E50").Value = "'  +0.63%  "
$ws.Range("D51").Value = "'6.936"
$ws.Range("E51").Value = "'  -1.16%  "
